# Threat Alert Report refresh - 2026-01-11 01:17
# Updates fare/threat figures for existing rows, reclassifies row 7 from
# "MEDIUM THREAT - MONITOR" down to "LOW THREAT", appends a new data row
# (row 9) for the 05-FEB-26 SM-329 flight, and narrows column J now that
# the longer "MEDIUM THREAT - MONITOR" text no longer appears there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a value into a "Date"-formatted column-A cell without
# Excel's autodetect silently turning the text into a date serial. We
# flip to Text format, assign the literal, then restore the original
# formatting (borders/fill/font/number-format) by pasting it in from a
# sibling cell that already carries the correct style - without then
# touching .Value again, which would re-trigger date parsing.
# ---------------------------------------------------------------------

# Row 2
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "15-JAN-26"
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("D2").Value = 12336
$ws.Range("E2").Value = 13395
$ws.Range("F2").Value = -1059

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "15-JAN-26"
$ws.Range("B3").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("C3").Value = "Air Arabia Egypt E5-585"
$ws.Range("D3").Value = 13072
$ws.Range("E3").Value = 13395
$ws.Range("F3").Value = -323

# Row 4
$ws.Range("D4").Value = 13244
$ws.Range("E4").Value = 13395
$ws.Range("F4").Value = -151

# Row 5
$ws.Range("D5").Value = 9060
$ws.Range("E5").Value = 9614
$ws.Range("F5").Value = -554

# Row 6
$ws.Range("D6").Value = 9350
$ws.Range("E6").Value = 9614
$ws.Range("F6").Value = -264

# Row 7
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "29-JAN-26"
$ws.Range("B7").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("C7").Value = "Air Arabia Egypt E5-585"
$ws.Range("D7").Value = 7200
$ws.Range("E7").Value = 7535
$ws.Range("F7").Value = -335
# Reclassify impact from MEDIUM THREAT - MONITOR to LOW THREAT, reusing
# the same formatting already used by the other LOW THREAT cells.
$ws.Range("J2").Copy()
$ws.Range("J7").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("J7").Value = "LOW THREAT"

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "31-JAN-26"
$ws.Range("B8").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("D8").Value = 8064
$ws.Range("E8").Value = 8455
$ws.Range("F8").Value = -391

# Row 9 (new) - clone formatting from row 8, then populate values
$ws.Range("A8:K8").Copy()
$ws.Range("A9:K9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "05-FEB-26"
$ws.Range("B9").Value = "SM-329"
$ws.Range("C9").Value = "Air Arabia Egypt E5-585"
$ws.Range("D9").Value = 7200
$ws.Range("E9").Value = 7535
$ws.Range("F9").Value = -335
$ws.Range("G9").Value = 30
$ws.Range("H9").Value = 30
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = "LOW THREAT"
$ws.Range("K9").Value = "EGP"
# Restore A9's formatting (copied row 8 already has the right style for
# every other cell); do this last so .Value is not touched afterwards.
$ws.Range("B9").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Column J is narrower now that it no longer needs to fit
# "MEDIUM THREAT - MONITOR" (OOXML width 25 -> 12; ColumnWidth runs
# ~0.83 narrower than the stored column width).
$ws.Columns.Item(10).ColumnWidth = 11.17
